# v1.2 update owner status
# LH_WF_PUBLISHVIDEO_REVIEWS workbook: update Owner Status (col H) values on the
# review sheet, and log the change in the VERSION-HISTORY sheet.

$wb = $excel.ActiveWorkbook
$wsReview   = $wb.Worksheets.Item("LH_WF_PUBLISHVIDEO_REVIEW")
$wsHistory  = $wb.Worksheets.Item("VERSION-HISTORY")

# --- LH_WF_PUBLISHVIDEO_REVIEW: update "Owner Status" (column H) --------------
$wsReview.Range("H2").Value = "closed"
$wsReview.Range("H4").Value = "closed"

# H3 previously had no fill (alternating-row style). Once it carries a real
# status it picks up the same light-blue fill used by the other rows in the
# striped table (matches the formatting already on G2/H2/H4 etc.) - copy just
# the formatting over from H2 so the theme-based fill is reused as-is, then
# set its own value.
$wsReview.Range("H2").Copy()
$wsReview.Range("H3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$wsReview.Range("H3").Value = "not applicable"

# The "select from list" validation on column H previously skipped H3 (it had
# no validation at all). Re-apply it across the full contiguous H2:H12 range.
$hRange = $wsReview.Range("H2:H12")
$hRange.Validation.Delete()
$hRange.Validation.Add(3, 1, 3, '"open,in progress,closed,not applicable"')
$hRange.Validation.ErrorTitle = "select from list"

# Leave the selection where the editor ended up.
$wsReview.Activate() | Out-Null
$wsReview.Range("C11").Select() | Out-Null

# --- VERSION-HISTORY: log the new version entry -------------------------------
$wsHistory.Range("A3").Value = "v1.2"
$wsHistory.Range("B3").Value = "eman"
$wsHistory.Range("C3").Value = "updtae owner status"
$wsHistory.Range("D3").Value2 = 45776

$wsHistory.Activate() | Out-Null
$wsHistory.Range("D12").Select() | Out-Null
